# "First batch of fixes after the feedback" - add the new feedback-related
# log entries (rows 21-23) to the hour log, matching the formatting of the
# rows already filled in above them, then refresh the total-hours formula's
# cached result and move the on-screen selection down to where the new rows
# were typed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 21: "Feedback" -----------------------------------------------
$ws.Range("A21").Value = "Feedback"
$ws.Range("B21").Value = 0
$ws.Range("C21").Value = "2025-05-19"
$ws.Range("D21").Value = "I just got my feedback on the assessment and i got some good improvements to work on. Since i handed in the first version of the assessment, i have worked on a big project where i made the API and the fetch functionaily, on a greater level then this assessment. So i have learned alot of new things that i can improve."

# Row 21's Subject cell is centered (like the "Done" row above it, row 17),
# unlike the other Subject cells which are left aligned.
$ws.Range("A21").HorizontalAlignment = -4108  # xlCenter
$ws.Range("A21").VerticalAlignment = -4108    # xlCenter
$ws.Range("D21").WrapText = $true
$ws.Range("D21").VerticalAlignment = -4108    # xlCenter

# --- Row 22: "Analyzed the feedback" -----------------------------------
$ws.Range("A22").Value = "Analyzed the feedback"
$ws.Range("B22").Value = 0
$ws.Range("C22").Value = "2025-05-19"
$ws.Range("D22").Value = "I took a while to read through the feedback and made a clear scheme of what to improve. Im going to improve them one at a time. About the feedback about the search functionality, i got this feedback: ""Search facilities allows me to search on each field but that is not what the assingment states."" But in the assessment, it says ""I want to be able to search facilities by the facility name, tag name, or location city, or any combination of those in a single API call."" I misunderstood this part whoopsie. Its a pretty easy fix luckily. The tags are the cause of the most issue so thats where the most time will be spent probably."

$ws.Range("A22").VerticalAlignment = -4108    # xlCenter
$ws.Range("D22").WrapText = $true
$ws.Range("D22").VerticalAlignment = -4108    # xlCenter

# --- Row 23: "Fixed the first couple issues" ---------------------------
$ws.Range("A23").Value = "Fixed the first couple issues"
$ws.Range("B23").Value = 2
$ws.Range("C23").Value = "2025-05-19"
$ws.Range("D23").Value = "I started with some of the easier fixes like the unnessecery search filters, made the tagets return as an array and fixed the database port issue. I also added documentation to the router page. Im going to add example bodys in the postman collection aswell to make it even more clear. Next up will probably be all the tags issues."

$ws.Range("A23").VerticalAlignment = -4108    # xlCenter
$ws.Range("D23").WrapText = $true
$ws.Range("D23").VerticalAlignment = -4108    # xlCenter

# --- Row heights: the newly-filled rows grew to fit their wrapped text,
# and the still-empty row 20 above them picked up the sheet's natural
# (non-custom) row height too. ------------------------------------------
$ws.Rows(20).RowHeight = 18
$ws.Rows(21).RowHeight = 38.25
$ws.Rows(22).RowHeight = 63.75
$ws.Rows(23).RowHeight = 38.25

# --- Move the selection to where editing left off -----------------------
$ws.Range("C19").Select()
